# Generate Report for Handoff
# Update Priority ("low" -> "ht") and refresh the Latest Handoff/HO-Xliff-Generate
# datetimes for the rows that were just re-handed-off.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7, Priority column E, Latest Handoff Datetime column H
foreach ($r in 4..7) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-10-21 05:20:09"
}

# de-de sheet: rows 4-7, Priority column E, Latest Handoff Datetime column H
foreach ($r in 4..7) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-10-21 05:20:21"
}

# Overview sheet: rows 4-7, Latest HO Xliff Generate Date column G shares the
# same underlying value as de-de!H4:H7.
foreach ($r in 4..7) {
    $overview.Range("G$r").Value = "2016-10-21 05:20:21"
}
